$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.052.96'
$ws.Range("E2").Value = '  +11.27%  '
$ws.Range("D3").Value = '1.814.85'
$ws.Range("E3").Value = '  +8.10%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '228.33'
$ws.Range("E5").Value = '  +3.71%  '
$ws.Range("E6").Value = '  +2.22%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '31.35'
$ws.Range("E8").Value = '  +4.14%  '
$ws.Range("D9").Value = '46.91'
$ws.Range("E9").Value = '  +5.89%  '
$ws.Range("E10").Value = '  +6.24%  '
$ws.Range("D11").Value = '0.0666'
$ws.Range("E11").Value = '  +4.32%  '
$ws.Range("E12").Value = '  +2.35%  '
$ws.Range("D13").Value = '2.076.31'
$ws.Range("E13").Value = '  +8.04%  '
$ws.Range("D14").Value = '1.822.30'
$ws.Range("E14").Value = '  +8.38%  '
$ws.Range("D15").Value = '0.640'
$ws.Range("E15").Value = '  +3.85%  '
$ws.Range("D16").Value = '34.061.74'
$ws.Range("E16").Value = '  +11.21%  '
$ws.Range("D17").Value = '10.26'
$ws.Range("D18").Value = '4.25'
$ws.Range("E18").Value = '  +6.52%  '
$ws.Range("D19").Value = '69.38'
$ws.Range("E19").Value = '  +4.29%  '
$ws.Range("D20").Value = '257.58'
$ws.Range("E20").Value = '  +4.84%  '
$ws.Range("D21").Value = '0.0₃0748'
$ws.Range("E21").Value = '  +3.27%  '
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("E23").Value = '  +4.11%  '
$ws.Range("D24").Value = '4.33'
$ws.Range("E24").Value = '  +0.99%  '
$ws.Range("E25").Value = '  +2.07%  '
$ws.Range("D26").Value = '158.38'
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").Value = '16.54'
$ws.Range("E27").Value = '  +3.77%  '
$ws.Range("D28").Value = '7.12'
$ws.Range("E28").Value = '  +6.09%  '
$ws.Range("D29").Value = '0.114'
$ws.Range("E29").Value = '  +1.46%  '
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("E31").Value = '  +10.40%  '
$ws.Range("E32").Value = '  +3.03%  '
$ws.Range("E33").Value = '  +4.58%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = '3.51'
$ws.Range("E34").Value = '  +6.29%  '
$ws.Range("B35").Value = 'MinaProtocolToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina'
$ws.Range("D35").Value = '1.55'
$ws.Range("E35").Value = '  +276.32%  '
$ws.Range("D36").Value = '1.537.90'
$ws.Range("E36").Value = '  +1.64%  '
$ws.Range("E37").Value = '  +2.19%  '
$ws.Range("E38").Value = '  +4.68%  '
$ws.Range("E39").Value = '  +0.82%  '
$ws.Range("E40").Value = '  +4.64%  '
$ws.Range("D41").Value = '0.621'
$ws.Range("E41").Value = '  +4.14%  '
$ws.Range("D42").Value = '2.81'
$ws.Range("E42").Value = '  +3.51%  '
$ws.Range("E43").Value = '  +1.66%  '
$ws.Range("D44").Value = '0.909'
$ws.Range("E44").Value = '  +8.04%  '
$ws.Range("D45").Value = '2.15'
$ws.Range("E45").Value = '  +7.73%  '
$ws.Range("E46").Value = '  +4.02%  '
$ws.Range("E47").Value = '  +4.82%  '
$ws.Range("D48").Value = '1.974.08'
$ws.Range("E48").Value = '  +8.50%  '
$ws.Range("D49").Value = '5.72'
$ws.Range("E49").Value = '  +1.78%  '
$ws.Range("D50").Value = '0.999'
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("D51").Value = '52.74'
$ws.Range("E51").Value = '  +1.60%  '
